$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values for rows 2-17 based on the diff (Natmi following Dr Hou advice)
$updates = @(
    @{ Row=2; E="3"; G="19.36022366666667"; H="58.080671"; I="0.005884129141485179"; J="0.005884129141485179"; K="3"; M="153.5290173333333"; N="460.587052"; O="0.3172206968818489"; P="0.317220696881849"; Q="2972.356114896877"; R="26751.20503407189"; S="0.001866567546804724"; T="0.001866567546804724" },
    @{ Row=3; E="3"; G="19.36022366666667"; H="58.080671"; I="0.005884129141485179"; J="0.005884129141485179"; K="3"; M="168.7997026666667"; N="506.3991080000001"; O="0.3487728915577651"; P="0.3487728915577651"; Q="3267.999998493497"; R="29411.99998644147"; S="0.002052224734975096"; T="0.002052224734975096" },
    @{ Row=4; E="3"; G="19.36022366666667"; H="58.080671"; I="0.005884129141485179"; J="0.005884129141485179"; K="3"; M="68.09032333333333"; N="204.27097"; O="0.1406878008722904"; P="0.1406878008722904"; Q="1318.243889268985"; R="11864.19500342087"; S="0.0008278251889641077"; T="0.0008278251889641079" },
    @{ Row=5; E="3"; G="19.36022366666667"; H="58.080671"; I="0.005884129141485179"; J="0.005884129141485179"; K="3"; M="93.562673"; N="280.688019"; O="0.1933186106880956"; P="0.1933186106880956"; Q="1811.394276131194"; R="16302.54848518075"; S="0.001137511670741252"; T="0.001137511670741252" },
    @{ Row=6; E="3"; G="3161.845459"; H="9485.536377"; I="0.9609758299542277"; J="0.9609758299542278"; K="3"; M="153.5290173333333"; N="460.587052"; O="0.3172206968818489"; P="0.317220696881849"; Q="485435.0262801322"; R="4368915.23652119"; S="0.3048414224646933"; T="0.3048414224646933" },
    @{ Row=7; E="3"; G="3161.845459"; H="9485.536377"; I="0.9609758299542277"; J="0.9609758299542278"; K="3"; M="168.7997026666667"; N="506.3991080000001"; O="0.3487728915577651"; P="0.3487728915577651"; Q="533718.5733571503"; R="4803467.160214352"; S="0.3351623189302592"; T="0.3351623189302592" },
    @{ Row=8; E="3"; G="3161.845459"; H="9485.536377"; I="0.9609758299542277"; J="0.9609758299542278"; K="3"; M="68.09032333333333"; N="204.27097"; O="0.1406878008722904"; P="0.1406878008722904"; Q="215291.0796333418"; R="1937619.716700076"; S="0.1351975762076844"; T="0.1351975762076844" },
    @{ Row=9; E="3"; G="3161.845459"; H="9485.536377"; I="0.9609758299542277"; J="0.9609758299542278"; K="3"; M="93.562673"; N="280.688019"; O="0.1933186106880956"; P="0.1933186106880956"; Q="295830.7127569519"; R="2662476.414812567"; S="0.1857745123515909"; T="0.185774512351591" },
    @{ Row=10; E="3"; G="2.055785333333333"; H="6.167356"; I="0.0006248123263850286"; J="0.0006248123263850286"; K="3"; M="153.5290173333333"; N="460.587052"; O="0.3172206968818489"; P="0.317220696881849"; Q="315.6227020749457"; R="2840.604318674512"; S="0.000198203401596228"; T="0.000198203401596228" },
    @{ Row=11; E="3"; G="2.055785333333333"; H="6.167356"; I="0.0006248123263850286"; J="0.0006248123263850286"; K="3"; M="168.7997026666667"; N="506.3991080000001"; O="0.3487728915577651"; P="0.3487728915577651"; Q="347.0159530131609"; R="3123.143577118448"; S="0.0002179176017542405"; T="0.0002179176017542405" },
    @{ Row=12; E="3"; G="2.055785333333333"; H="6.167356"; I="0.0006248123263850286"; J="0.0006248123263850286"; K="3"; M="68.09032333333333"; N="204.27097"; O="0.1406878008722904"; P="0.1406878008722904"; Q="139.9790880505911"; R="1259.81179245532"; S="8.79034721570094E-05"; T="8.790347215700943E-05" },
    @{ Row=13; E="3"; G="2.055785333333333"; H="6.167356"; I="0.0006248123263850286"; J="0.0006248123263850286"; K="3"; M="93.562673"; N="280.688019"; O="0.1933186106880956"; P="0.1933186106880956"; Q="192.3447709008626"; R="1731.102938107764"; S="0.0001207878508775507"; T="0.0001207878508775507" },
    @{ Row=14; E="3"; G="106.9830526666667"; H="320.949158"; I="0.03251522857790212"; J="0.03251522857790212"; K="3"; M="153.5290173333333"; N="460.587052"; O="0.3172206968818489"; P="0.317220696881849"; Q="16425.00294723358"; R="147825.0265251022"; S="0.01031450346875472"; T="0.01031450346875472" },
    @{ Row=15; E="3"; G="106.9830526666667"; H="320.949158"; I="0.03251522857790212"; J="0.03251522857790212"; K="3"; M="168.7997026666667"; N="506.3991080000001"; O="0.3487728915577651"; P="0.3487728915577651"; Q="18058.70748050568"; R="162528.3673245511"; S="0.0113404302907766"; T="0.0113404302907766" },
    @{ Row=16; E="3"; G="106.9830526666667"; H="320.949158"; I="0.03251522857790212"; J="0.03251522857790212"; K="3"; M="68.09032333333333"; N="204.27097"; O="0.1406878008722904"; P="0.1406878008722904"; Q="7284.510647260362"; R="65560.59582534326"; S="0.004574496003484899"; T="0.0045744960034849" },
    @{ Row=17; E="3"; G="106.9830526666667"; H="320.949158"; I="0.03251522857790212"; J="0.03251522857790212"; K="3"; M="93.562673"; N="280.688019"; O="0.1933186106880956"; P="0.1933186106880956"; Q="10009.62037319311"; R="90086.583358738"; S="0.006285798814885902"; T="0.006285798814885902" }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 5).Value = [double]$u.E
    $ws.Cells.Item($r, 7).Value = [double]$u.G
    $ws.Cells.Item($r, 8).Value = [double]$u.H
    $ws.Cells.Item($r, 9).Value = [double]$u.I
    $ws.Cells.Item($r, 10).Value = [double]$u.J
    $ws.Cells.Item($r, 11).Value = [double]$u.K
    $ws.Cells.Item($r, 13).Value = [double]$u.M
    $ws.Cells.Item($r, 14).Value = [double]$u.N
    $ws.Cells.Item($r, 15).Value = [double]$u.O
    $ws.Cells.Item($r, 16).Value = [double]$u.P
    $ws.Cells.Item($r, 17).Value = [double]$u.Q
    $ws.Cells.Item($r, 18).Value = [double]$u.R
    $ws.Cells.Item($r, 19).Value = [double]$u.S
    $ws.Cells.Item($r, 20).Value = [double]$u.T
}
